$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216; existing rows 216-245 shift down to 217-246.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new weekly price record.
$ws.Range("A216").Value = 10
$ws.Range("B216").Value = "Vega Modelo de Temuco"
$ws.Range("C216").Value = "La Araucanía"
$ws.Range("D216").Value = 44449
$ws.Range("E216").Value = 9
$ws.Range("F216").Value = 100114014
$ws.Range("G216").Value = "Betarraga"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 45
$ws.Range("K216").Value = 8000
$ws.Range("L216").Value = 8000
$ws.Range("M216").Value = 8000
$ws.Range("N216").Value = "$/docena de paquetes"
$ws.Range("O216").Value = "Provincia de Cautín"
$ws.Range("P216").Value = 667
$ws.Range("Q216").Value = 12
$ws.Range("R216").Value = "Hortaliza"
